$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Memory table: row 52 gains a "BankID 4 Bit" column (C52) and the
#     read_iram instruction is renamed to read_ram (E52) ---
$ws.Range("C52").Value = "BankID 4 Bit"
$ws.Range("E52").Value = "read_ram"

# --- New row 53: 0x33 / VRAM Write ---
$ws.Range("B53").Value = "0x33"
$ws.Range("C53").Value = "Addr"
$ws.Range("D53").Value = "Data"
$ws.Range("E53").Value = "VRAM Write"
$ws.Range("G53").Value = "0x33"

# --- New row 54: 0x34 / MMIO Write ---
$ws.Range("B54").Value = "0x34"
$ws.Range("C54").Value = "Addr"
$ws.Range("D54").Value = "Data"
$ws.Range("E54").Value = "MMIO Write"
$ws.Range("G54").Value = "0x34"

# --- View state: scroll back to the top and move the active selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select()

# --- Minor column width tweaks (B:C, E, F) ---
$ws.Columns.Item(2).ColumnWidth = 16.3
$ws.Columns.Item(3).ColumnWidth = 16.3
$ws.Columns.Item(5).ColumnWidth = 42.3
$ws.Columns.Item(6).ColumnWidth = 16.3

Write-Output "DebugDataTransfer updated"
